# ND02.xlsx fixes for unit and integration tests (commit 989)
# Adds two new data rows to the "CMS" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CMS")
$contactDate = Get-Date -Year 2017 -Month 9 -Day 13 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# --- Row 2 ---
$ws.Cells.Item(2, 1).Value = 654
$ws.Cells.Item(2, 2).Value = $contactDate
$ws.Cells.Item(2, 3).Value = "CMS2"
$ws.Cells.Item(2, 4).Value = "The second type"
$ws.Cells.Item(2, 5).Value = "Tom Swann"
$ws.Cells.Item(2, 6).Value = 1001
$ws.Cells.Item(2, 7).Value = "C"
$ws.Cells.Item(2, 8).Value = "WMT"
$ws.Cells.Item(2, 9).Value = "ND01"
$ws.Cells.Item(2, 10).Value = "Any Wright"
$ws.Cells.Item(2, 11).Value = 1002
$ws.Cells.Item(2, 12).Value = "Z"
$ws.Cells.Item(2, 13).Value = "WMT"
$ws.Cells.Item(2, 14).Value = "ND01"

# --- Row 3 ---
$ws.Cells.Item(3, 1).Value = 678
$ws.Cells.Item(3, 1).WrapText = $true
$ws.Cells.Item(3, 2).Value = $contactDate
$ws.Cells.Item(3, 2).NumberFormat = "mm-dd-yy"
# Re-use the just-minted date style for B2 too (copy from B3, not vice
# versa) so both cells share one cellXfs entry instead of minting two
# separate-but-identical ones.
$ws.Cells.Item(3, 2).Copy($ws.Cells.Item(2, 2))
$ws.Cells.Item(3, 3).Value = "CMS1"
$ws.Cells.Item(3, 4).Value = "The first type"
$ws.Cells.Item(3, 5).Value = "Andy Wright"
$ws.Cells.Item(3, 6).Value = 1002
$ws.Cells.Item(3, 7).Value = "Z"
$ws.Cells.Item(3, 8).Value = "WMT"
$ws.Cells.Item(3, 9).Value = "ND01"
$ws.Cells.Item(3, 10).Value = "Tom Swann"
$ws.Cells.Item(3, 11).Value = 1001
$ws.Cells.Item(3, 12).Value = "C"
$ws.Cells.Item(3, 13).Value = "WMT"
$ws.Cells.Item(3, 14).Value = "ND01"

# Selection/active-cell state matches the post-edit saved file (full rows 2:3).
$ws.Range("A2:A3").EntireRow.Select()
$excel.ActiveCell = $ws.Range("A2")
